$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 03:16"

# Row 7 - Espana
$ws.Range("B7").Value = 28768
$ws.Range("C7").Value = 3272
$ws.Range("D7").Value = 2575
$ws.Range("E7").Value = 24421
$ws.Range("G7").Value = 391
$ws.Range("H7").Value = 1772

# Row 9 - Iran
$ws.Range("D9").Value = 7913
$ws.Range("E9").Value = 12040

# Row 15 - Austria
$ws.Range("B15").Value = 3582
$ws.Range("C15").Value = 590
$ws.Range("E15").Value = 3557

# Row 17 - Noruega
$ws.Range("B17").Value = 2385
$ws.Range("C17").Value = 221
$ws.Range("E17").Value = 2372

# Rows 19-23: Australia's updated case count overtakes Portugal, Brasil,
# Canada and Dinamarca, so it is re-sorted into the table right after
# Suecia (row 18) and the other four countries shift down one row,
# keeping their own (unchanged) figures.

# Row 19 - now Australia (its own new, updated values)
$ws.Range("A19").Value = "Australia"
$ws.Range("B19").Value = 1609
$ws.Range("C19").Value = 537
$ws.Range("D19").Value = 88
$ws.Range("E19").Value = 1514
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 7

# Row 20 - now Portugal (its own, unchanged values, shifted down from row 19)
$ws.Range("A20").Value = "Portugal"
$ws.Range("B20").Value = 1600
$ws.Range("C20").Value = 320
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 1581
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 14

# Row 21 - now Brasil (its own, unchanged values, shifted down from row 20)
$ws.Range("A21").Value = "Brasil"
$ws.Range("B21").Value = 1546
$ws.Range("C21").Value = 368
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 1519
$ws.Range("F21").Value = 18
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 25

# Row 22 - now Canada (its own, unchanged values, shifted down from row 21)
$ws.Range("A22").Value = "Canada"
$ws.Range("B22").Value = 1470
$ws.Range("C22").Value = 142
$ws.Range("D22").Value = 14
$ws.Range("E22").Value = 1436
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 20

# Row 23 - now Dinamarca (its own, unchanged values, shifted down from row 22)
$ws.Range("A23").Value = "Dinamarca"
$ws.Range("B23").Value = 1395
$ws.Range("C23").Value = 69
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1381
$ws.Range("F23").Value = 42
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 13

# Row 27 - Japon
$ws.Range("B27").Value = 1101
$ws.Range("C27").Value = 47
$ws.Range("E27").Value = 825
$ws.Range("F27").Value = 49
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 41
